# REPORTGEN-1102: part 1, added and removed counts missing when no previous snapshot selected
#
# The "evolution" (added/removed vulnerabilities) tables and the quality-standards
# evolution table were not being flagged as such, so RepGen could not compute
# added/removed counts when no previous snapshot was selected. This appends the
# ",EVOLUTION=true" parameter to the relevant RepGen table markers on the
# Summary sheet and on each "Ax-2017" overview sheet.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Summary").Range("B14").Value = "RepGen:TABLE;QUALITY_STANDARDS_EVOLUTION;STD=OWASP-2017,HEADER=NO,EVOLUTION=true"

$wb.Worksheets.Item("A1-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A1-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A2-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A2-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A3-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A3-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A4-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A4-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A5-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A5-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A6-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A6-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A7-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A7-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A8-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A8-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A9-2017").Range("A3").Value  = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A9-2017,DESC=true,HEADER=NO,EVOLUTION=true"
$wb.Worksheets.Item("A10-2017").Range("A3").Value = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=A10-2017,DESC=true,HEADER=NO,EVOLUTION=true"
